$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the daily conversion note text on Hoja1!A1 with the new rates.
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.38 = 54496.79 pesos`n✅ 54496.79 pesos = 13.38 = 984.58 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the "tasas" sheet rate figures.
$ws2.Range("N10").Value = 74.72
$ws2.Range("O10").Value = 4072
$ws2.Range("N12").Value = 4073.8
$ws2.Range("O12").Value = 73.59999999999999
